$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values (B11, B12, B13 re-sorted)
$ws.Range("B11").Value = 5.57
$ws.Range("B12").Value = 5.62
$ws.Range("B13").Value = 5.66

# B11 gains underline font + "0.00" number format (new combined style)
$ws.Range("B11").NumberFormat = "0.00"
$ws.Range("B11").Font.Underline = $true

# B13 loses its underline formatting, back to default/general
$ws.Range("B13").Font.Underline = $false

# Update selection to B11
$ws.Range("B11").Select()
